$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking price strings
# (e.g. "195.08") are stored as text, matching the inlineStr cells in the source file.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '76.107.72'
$ws.Range("E2").Value = '  +0.90%  '

$ws.Range("D3").Value = '2.859.52'
$ws.Range("E3").Value = '  +7.38%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = '195.08'
$ws.Range("E5").Value = '  +3.96%  '

$ws.Range("D6").Value = '596.54'
$ws.Range("E6").Value = '  +1.59%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '0.553'
$ws.Range("E8").Value = '  +3.60%  '

$ws.Range("D9").Value = '0.192'
$ws.Range("E9").Value = '  -1.11%  '

$ws.Range("D10").Value = '2.853.66'
$ws.Range("E10").Value = '  +7.18%  '

$ws.Range("D11").Value = '0.391'
$ws.Range("E11").Value = '  +9.60%  '

$ws.Range("E12").Value = '  -2.03%  '

$ws.Range("D13").Value = '4.90'
$ws.Range("E13").Value = '  +3.72%  '

$ws.Range("D14").Value = '3.386.33'
$ws.Range("E14").Value = '  +7.37%  '

$ws.Range("D15").Value = '76.036.18'
$ws.Range("E15").Value = '  +1.05%  '

$ws.Range("D16").Value = '27.40'
$ws.Range("E16").Value = '  +3.20%  '

$ws.Range("D17").Value = '0.0000188'
$ws.Range("E17").Value = '  -0.07%  '

$ws.Range("D18").Value = '2.875.12'
$ws.Range("E18").Value = '  +7.59%  '

$ws.Range("D19").Value = '9.01'
$ws.Range("E19").Value = '  -1.57%  '

$ws.Range("D20").Value = '12.48'
$ws.Range("E20").Value = '  +4.53%  '

$ws.Range("D21").Value = '380.34'
$ws.Range("E21").Value = '  +2.42%  '

$ws.Range("D22").Value = '2.32'
$ws.Range("E22").Value = '  +2.11%  '

$ws.Range("D23").Value = '4.12'
$ws.Range("E23").Value = '  +0.97%  '

$ws.Range("D24").Value = '71.62'
$ws.Range("E24").Value = '  +2.60%  '

$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("D26").Value = '3.029.40'
$ws.Range("E26").Value = '  +8.07%  '

$ws.Range("D27").Value = '4.19'
$ws.Range("E27").Value = '  +1.03%  '

$ws.Range("D28").Value = '9.71'
$ws.Range("E28").Value = '  +4.45%  '

$ws.Range("D29").Value = '0.0000104'
$ws.Range("E29").Value = '  +10.01%  '

$ws.Range("E30").Value = '  +0.33%  '

$ws.Range("D31").Value = '1.40'
$ws.Range("E31").Value = '  -1.57%  '

$ws.Range("D32").Value = '508.10'
$ws.Range("E32").Value = '  -2.04%  '

$ws.Range("D33").Value = '7.69'
$ws.Range("E33").Value = '  -0.09%  '

$ws.Range("D34").Value = '1.80'
$ws.Range("E34").Value = '  +2.97%  '

$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("D36").Value = '167.34'
$ws.Range("E36").Value = '  +2.41%  '

$ws.Range("D37").Value = '19.94'
$ws.Range("E37").Value = '  +4.02%  '

$ws.Range("E38").Value = '  -1.25%  '

$ws.Range("D39").Value = '19.54'
$ws.Range("E39").Value = '  +0.84%  '

$ws.Range("D40").Value = '183.90'
$ws.Range("E40").Value = '  +8.51%  '

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("D42").Value = '0.344'
$ws.Range("E42").Value = '  +4.52%  '

$ws.Range("D43").Value = '5.03'
$ws.Range("E43").Value = '  +1.05%  '

$ws.Range("D44").Value = '1.67'
$ws.Range("E44").Value = '  -1.45%  '

$ws.Range("D45").Value = '0.0919'
$ws.Range("E45").Value = '  +8.76%  '

$ws.Range("D46").Value = '1.22'
$ws.Range("E46").Value = '  +2.18%  '

$ws.Range("D47").Value = '40.18'
$ws.Range("E47").Value = '  +2.78%  '

$ws.Range("D48").Value = '2.34'
$ws.Range("E48").Value = '  -1.19%  '

$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '0.577'
$ws.Range("E49").Value = '  +8.41%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '0.682'
$ws.Range("E50").Value = '  +15.58%  '

$ws.Range("D51").Value = '3.73'
$ws.Range("E51").Value = '  +2.52%  '

# Reset style index back to the default "Normal" style so no stray cell-level
# style attribute is left behind (keeps styles.xml / cell "s" refs unchanged).
$ws.Range("D2:D51").Style = "Normal"
